$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 9312.25
$ws.Range("I86").Value = 3250
$ws.Range("J86").Value = 11333
$ws.Range("K86").Value = 3250
$ws.Range("L86").Value = 11333
$ws.Range("M86").Value = -2127
$ws.Range("N86").Value = -13579

$ws.Range("H89").Value = 9312.25
$ws.Range("I89").Value = 3250
$ws.Range("J89").Value = 11333
$ws.Range("K89").Value = 16250
$ws.Range("L89").Value = 56665
$ws.Range("M89").Value = -10634
$ws.Range("N89").Value = -67897

$ws.Range("H100").Value = 1694.6842
$ws.Range("I100").Value = 1885.5
$ws.Range("J100").Value = 1160.4
$ws.Range("K100").Value = 1885.5
$ws.Range("L100").Value = 1160.4
$ws.Range("M100").Value = -1344.5
$ws.Range("N100").Value = -2242.4

$ws.Range("H113").Value = 11168.111
$ws.Range("I113").Value = 12339.25
$ws.Range("J113").Value = 1799
$ws.Range("K113").Value = 12339.25
$ws.Range("L113").Value = 1799
$ws.Range("M113").Value = -9085.25
$ws.Range("N113").Value = -8307

$ws.Range("H127").Value = 920.06665
$ws.Range("I127").Value = 829.53845
$ws.Range("J127").Value = 1508.5
$ws.Range("K127").Value = 2488.61535
$ws.Range("L127").Value = 4525.5
$ws.Range("M127").Value = 2471.38465
$ws.Range("N127").Value = -14445.5

$ws.Range("H129").Value = 58685.285
$ws.Range("I129").Value = 50493
$ws.Range("K129").Value = 151479
$ws.Range("M129").Value = -146479

$ws.Range("H137").Value = 68212.05
$ws.Range("I137").Value = 8076.6924
$ws.Range("J137").Value = 97166.11
$ws.Range("K137").Value = 24230.0772
$ws.Range("L137").Value = 291498.33
$ws.Range("M137").Value = -21680.0772
$ws.Range("N137").Value = -296598.33

$ws.Range("H138").Value = 4502.2964
$ws.Range("I138").Value = 6596.2856
$ws.Range("J138").Value = 4190.4253
$ws.Range("K138").Value = 19788.8568
$ws.Range("L138").Value = 12571.2759
$ws.Range("M138").Value = -14648.8568
$ws.Range("N138").Value = -22851.2759

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 882.65
$ws.Range("I2").Value = 823.8421
$ws.Range("J2").Value = 2000
$ws.Range("K2").Value = 823.8421
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = -710.8421
$ws.Range("N2").Value = -2226

$ws.Range("H61").Value = 1350674.5
$ws.Range("I61").Value = 4817.1333
$ws.Range("J61").Value = 4234654.5
$ws.Range("K61").Value = 4817.1333
$ws.Range("L61").Value = 4234654.5
$ws.Range("M61").Value = -4605.1333
$ws.Range("N61").Value = -4235078.5

$ws.Range("H97").Value = 13634.667
$ws.Range("I97").Value = 9603
$ws.Range("J97").Value = 17666.334
$ws.Range("K97").Value = 9603
$ws.Range("L97").Value = 17666.334
$ws.Range("M97").Value = -9107
$ws.Range("N97").Value = -18658.334

$ws.Range("H102").Value = 2008.7273
$ws.Range("I102").Value = 2324.5
$ws.Range("J102").Value = 1166.6666
$ws.Range("K102").Value = 2324.5
$ws.Range("L102").Value = 1166.6666
$ws.Range("M102").Value = -702.5
$ws.Range("N102").Value = -4410.6666

$ws.Range("H116").Value = 882.65
$ws.Range("I116").Value = 823.8421
$ws.Range("J116").Value = 2000
$ws.Range("K116").Value = 823.8421
$ws.Range("L116").Value = 2000
$ws.Range("M116").Value = 1470.1579
$ws.Range("N116").Value = -6588

$ws.Range("H136").Value = 1350674.5
$ws.Range("I136").Value = 4817.1333
$ws.Range("J136").Value = 4234654.5
$ws.Range("K136").Value = 14451.3999
$ws.Range("L136").Value = 12703963.5
$ws.Range("M136").Value = -11901.3999
$ws.Range("N136").Value = -12709063.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 882.65
$ws.Range("I3").Value = 823.8421
$ws.Range("J3").Value = 2000
$ws.Range("K3").Value = 823.8421
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = -709.8421
$ws.Range("N3").Value = -2228

$ws.Range("H94").Value = 4772.2856
$ws.Range("I94").Value = 3964
$ws.Range("J94").Value = 5850
$ws.Range("K94").Value = 3964
$ws.Range("L94").Value = 5850
$ws.Range("M94").Value = -3513
$ws.Range("N94").Value = -6752

$ws.Range("H134").Value = 8952.035
$ws.Range("I134").Value = 5709.125
$ws.Range("K134").Value = 17127.375
$ws.Range("M134").Value = -14592.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 16186.77
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 16186.77
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 16186.77
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -16776.77

$ws.Range("H34").Value = 16186.77
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 16186.77
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 16186.77
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -16590.77

$ws.Range("H99").Value = 8554.210999999999
$ws.Range("I99").Value = 2144.5
$ws.Range("J99").Value = 19542.285
$ws.Range("K99").Value = 2144.5
$ws.Range("L99").Value = 19542.285
$ws.Range("M99").Value = -646.5
$ws.Range("N99").Value = -22538.285

$ws.Range("H122").Value = 2733.4138
$ws.Range("I122").Value = 1107.6
$ws.Range("J122").Value = 6346.3335
$ws.Range("K122").Value = 3322.8
$ws.Range("L122").Value = 19039.0005
$ws.Range("M122").Value = -872.7999999999997
$ws.Range("N122").Value = -23939.0005

$ws.Range("H126").Value = 8554.210999999999
$ws.Range("I126").Value = 2144.5
$ws.Range("J126").Value = 19542.285
$ws.Range("K126").Value = 6433.5
$ws.Range("L126").Value = 58626.855
$ws.Range("M126").Value = -3963.5
$ws.Range("N126").Value = -63566.855

$ws.Range("H134").Value = 45462524
$ws.Range("I134").Value = 2370.6667
$ws.Range("J134").Value = 142877140
$ws.Range("K134").Value = 7112.000100000001
$ws.Range("L134").Value = 428631420
$ws.Range("M134").Value = -4577.000100000001
$ws.Range("N134").Value = -428636490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 3682.4595
$ws.Range("J68").Value = 3904.8125
$ws.Range("L68").Value = 11714.4375
$ws.Range("N68").Value = -13336.4375

$ws.Range("H71").Value = 3682.4595
$ws.Range("J71").Value = 3904.8125
$ws.Range("L71").Value = 35143.3125
$ws.Range("N71").Value = -43255.3125

$ws.Range("H122").Value = 12663251
$ws.Range("J122").Value = 3154362
$ws.Range("L122").Value = 28389258
$ws.Range("N122").Value = -28394158

$ws.Range("H131").Value = 1455.74
$ws.Range("J131").Value = 1498.7957
$ws.Range("L131").Value = 4496.3871
$ws.Range("N131").Value = -14576.3871

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H117").Value = 50000
$ws.Range("J117").Value = 50000
$ws.Range("L117").Value = 50000
$ws.Range("N117").Value = -56884

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4901.846
$ws.Range("J22").Value = 5256.1665
$ws.Range("L22").Value = 5256.1665
$ws.Range("N22").Value = -5846.1665

$ws.Range("H27").Value = 4901.846
$ws.Range("J27").Value = 5256.1665
$ws.Range("L27").Value = 5256.1665
$ws.Range("N27").Value = -5470.1665

$ws.Range("H93").Value = 27800.75
$ws.Range("I93").Value = 27800.75
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 27800.75
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -26552.75
$ws.Range("N93").ClearContents()

$ws.Range("H128").Value = 70000
$ws.Range("J128").Value = 70000
$ws.Range("L128").Value = 70000
$ws.Range("N128").Value = -79960

$ws.Range("H132").Value = 3495927
$ws.Range("I132").Value = 7240
$ws.Range("J132").Value = 8728957
$ws.Range("K132").Value = 21720
$ws.Range("L132").Value = 26186871
$ws.Range("M132").Value = -19190
$ws.Range("N132").Value = -26191931

$ws.Range("H136").Value = 1456639
$ws.Range("I136").Value = 25807
$ws.Range("J136").Value = 2887471
$ws.Range("K136").Value = 77421
$ws.Range("L136").Value = 8662413
$ws.Range("M136").Value = -74871
$ws.Range("N136").Value = -8667513

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1499
$ws.Range("I96").Value = 1499
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 1499
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -126
$ws.Range("N96").ClearContents()
